$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that needs to be bumped
# from 45180 (2023-09-11) to 45181 (2023-09-12) for every data row
# (rows 2 through 232).
$ws.Range("C2:C232").Value = 45181
